# Update Bitcoin price history through 2023-09-09
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the previously-last row (B1326) per upstream source recompute
$ws.Range("B1326").Value = 27895.54510920212

# Append rows 1327-1349 (dates 2023-08-18 .. 2023-09-09 with their prices)
$newDates = @(
  "2023-08-18",
  "2023-08-19",
  "2023-08-20",
  "2023-08-21",
  "2023-08-22",
  "2023-08-23",
  "2023-08-24",
  "2023-08-25",
  "2023-08-26",
  "2023-08-27",
  "2023-08-28",
  "2023-08-29",
  "2023-08-30",
  "2023-08-31",
  "2023-09-01",
  "2023-09-02",
  "2023-09-03",
  "2023-09-04",
  "2023-09-05",
  "2023-09-06",
  "2023-09-07",
  "2023-09-08",
  "2023-09-09"
)
$newValues = @(
  26040.81931205392,
  26115.91295705558,
  26137.35350149368,
  26086.51127135073,
  25826.03364121723,
  26641.12006805851,
  26021.08163211241,
  25949.31275796248,
  26025.2295197059,
  26112.34304465044,
  25975.58230041972,
  27852.18330418599,
  27251.96914888735,
  26179.07181383978,
  25581.65193106596,
  25845.21631827925,
  25991.65745373537,
  25871.5364941652,
  25706.43383971701,
  25643.30464824743,
  25871.77707088979,
  25880.11091714333,
  25857.72787564452
)

$startRow = 1327
for ($i = 0; $i -lt $newDates.Length; $i++) {
  $row = $startRow + $i
  $aCell = $ws.Range("A$row")
  # Use a text-producing formula so the date-like string is not
  # auto-converted to a date serial, then collapse it to a plain
  # static value (matches the original sheet's plain shared-string cells).
  $aCell.Formula = '="' + $newDates[$i] + '"'
  $aCell.Copy()
  $aCell.PasteSpecial(-4163)

  $ws.Range("B$row").Value = $newValues[$i]
}

$excel.CutCopyMode = 0
